# TalbertTso Performance Appraisal - apply commit edits
$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Extend the "student worker Tester group" paragraph with two new
#    sentences (inherits the Arial formatting of the run being extended).
# ---------------------------------------------------------------------------
$addition = "possibly assist in updating older application. " + `
  "Our team did not have any prior job description or interview materials " + `
  "at the ready to immediately start seeking or hiring student workers in " + `
  "this specific field."

$null = $d.Content.Find.Execute(
  "assist in any testing groups that might occur, and ",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "assist in any testing groups that might occur, and " + $addition,
  2)

# ---------------------------------------------------------------------------
# 2) The blank paragraph right after that table cell's text now carries an
#    explicit Arial rPr on its paragraph mark. Locate it via the paragraph
#    that follows the one we just edited.
# ---------------------------------------------------------------------------
$r2 = $d.Content
$null = $r2.Find.Execute("Our team did not have any prior job description or interview materials at the ready to immediately start seeking or hiring student workers in this specific field.")
$r2.Collapse(0)
$nextPara = $r2.Paragraphs(1)
$nextPara.Range.Font.Name = "Arial"

# ---------------------------------------------------------------------------
# 3) Drop the stray "lastRenderedPageBreak" rendering caches that no longer
#    correspond to the (now different) pagination. A same-text Find/Replace
#    rewrites the run and sheds the cached break marker.
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
  "Ensures own actions are consistent with the university",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "Ensures own actions are consistent with the university",
  2)

$null = $d.Content.Find.Execute(
  "List training and development opportunities participated in during this appraisal period",
  $true, $false, $false, $false, $false, $true, 1, $false,
  "List training and development opportunities participated in during this appraisal period",
  2)

# ---------------------------------------------------------------------------
# 4) Footer page-count field result cache: stale "3" -> "6".
# ---------------------------------------------------------------------------
$ftr = $d.Sections(1).Footers(1).Range
$null = $ftr.Find.Execute("3", $true, $false, $false, $false, $false, $true, 1, $false, "6", 2)
